$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Append the four new result rows (56-59) below the existing data.
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row=56; A=1; B=4; C=50; D=18; E=2; F=17078;  G=17302;  H=0; I=129532;  J=0; K=0; L=18.547377999999998; M="SPYlearner"; N="ExtraStates:1+EQ"; O=1; P="TeacherDFSM"; R="DFA_R50_peterson2.fsm" },
    @{ Row=57; A=1; B=4; C=50; D=18; E=2; F=293439; G=293816; H=0; I=2495242; J=0; K=0; L=28.324876;           M="SPYlearner"; N="ExtraStates:2+EQ"; O=2; P="TeacherDFSM"; R="DFA_R50_peterson2.fsm" },
    @{ Row=58; A=1; B=4; C=97; D=12; E=2; F=16498;  G=17106;  H=0; I=207606;  J=0; K=0; L=35.058551999999999; M="SPYlearner"; N="ExtraStates:1+EQ"; O=1; P="TeacherDFSM"; R="DFA_R97_sched4.fsm" },
    @{ Row=59; A=1; B=4; C=97; D=12; E=2; F=172412; G=173577; H=0; I=2297891; J=0; K=0; L=48.350580000000001; M="SPYlearner"; N="ExtraStates:2+EQ"; O=2; P="TeacherDFSM"; R="DFA_R97_sched4.fsm" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 18).Value = $r.R
}

# ---------------------------------------------------------------------------
# 2. Apply an AutoFilter over the full (now-grown) data range, filtering
#    column C ("States") down to the discrete value 97.
# ---------------------------------------------------------------------------
$fullRange = $ws.Range("A1:R59")
$fullRange.AutoFilter(3, @("97"), 7)

# ---------------------------------------------------------------------------
# 3. Register the hidden sheet-scoped _FilterDatabase defined name that Excel
#    writes whenever an AutoFilter is in effect.
# ---------------------------------------------------------------------------
$fdName = $ws.Names.Add("_xlnm._FilterDatabase", "=results!`$A`$1:`$R`$59")
$fdName.Visible = $false

# ---------------------------------------------------------------------------
# 4. Update the visible selection to match the post-filter cursor position.
# ---------------------------------------------------------------------------
$ws.Range("B39").Select()
